# CDS Data Validations script updates
# Replaces the "WebExcel" data-source column with an "ExDataExcel" column,
# pointing at the generated "_ExcelData.xlsx" file instead of "_WebData.xlsx".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D header: WebExcel -> ExDataExcel ---
$ws.Range("D1").Value = "ExDataExcel"

# --- Column D data rows: *_WebData.xlsx -> *_ExcelData.xlsx ---
$newFileName = "TC04_CDSValidation_by_ParticipantID - 5_ExcelData.xlsx"
$ws.Range("D2").Value = $newFileName
$ws.Range("D3").Value = $newFileName
$ws.Range("D4").Value = $newFileName

# --- Widen column C to fit the new, longer header/values ---
$ws.Columns.Item(3).ColumnWidth = 54.25

# --- Update the view: scroll down a row and select D2 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 2
$win.ScrollColumn = 1
$ws.Range("D2").Select()
